$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")


# ---- 展览 ----
$wsExhibit.Range("F2").Value = 12759
$wsExhibit.Range("F3").Value = 7127
$wsExhibit.Range("F10").Value = 1000
$wsExhibit.Range("F11").Value = 142
$wsExhibit.Range("F12").Value = 350
$wsExhibit.Range("F13").Value = 1006
$wsExhibit.Range("F14").Value = 2
$wsExhibit.Range("F16").Value = 1015
$wsExhibit.Range("F18").Value = 241
$wsExhibit.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202408/QupCAJLM1723713571010.jpeg"
$wsExhibit.Range("F22").Value = 306
$wsExhibit.Range("F24").Value = 142
$wsExhibit.Range("G24").Value = 88
$wsExhibit.Range("F26").Value = 5220
$wsExhibit.Range("F27").Value = 70
$wsExhibit.Range("F28").Value = 1421
$wsExhibit.Range("F29").Value = 306
$wsExhibit.Range("F30").Value = 1329
$wsExhibit.Range("F31").Value = 56
$wsExhibit.Range("F32").Value = 28
$wsExhibit.Range("F33").Value = 1352
$wsExhibit.Range("F36").Value = 591
$wsExhibit.Range("F38").Value = 3727

# ---- 演出 ----
$wsShow.Range("F4").Value = 3741
$wsShow.Range("F5").Value = 3741
$wsShow.Range("F8").Value = 52
$wsShow.Range("C10").Value = "杭州·《天空之城》动漫经典音乐作品演奏会（取消）"
$wsShow.Range("G10").Value = "不可售"

# ---- 本地生活 ----
$wsLocal.Range("F2").Value = 9265
$wsLocal.Range("F3").Value = 557
$wsLocal.Range("F4").Value = 1997

# ---- 全部类型 ----
$wsAll.Range("F2").Value = 9265
$wsAll.Range("F3").Value = 557
$wsAll.Range("F4").Value = 1997
$wsAll.Range("F5").Value = 12760
$wsAll.Range("F6").Value = 7127
$wsAll.Range("F8").Value = 3741
$wsAll.Range("F11").Value = 142
$wsAll.Range("F12").Value = 350
$wsAll.Range("F13").Value = 1006
$wsAll.Range("F14").Value = 2
$wsAll.Range("F16").Value = 1015
$wsAll.Range("F18").Value = 241
$wsAll.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202408/QupCAJLM1723713571010.jpeg"
$wsAll.Range("F22").Value = 306
$wsAll.Range("C27").Value = "杭州·2024首届COMIC GALAXY次元盛典"
$wsAll.Range("D27").Value = "长江南路336号 白马湖国际会展中心"
$wsAll.Range("E27").Value = "2024.09.15 09:30-09.17 17:30"
$wsAll.Range("F27").Value = 142
$wsAll.Range("G27").Value = 88
$wsAll.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=90433"
$wsAll.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202408/teoBMbzd1723019674766.png"
$wsAll.Range("F29").Value = 5220
$wsAll.Range("F30").Value = 70
$wsAll.Range("F31").Value = 1421
$wsAll.Range("F34").Value = 306
$wsAll.Range("F36").Value = 1329
$wsAll.Range("F37").Value = 56
$wsAll.Range("F38").Value = 1352
$wsAll.Range("F40").Value = 591
$wsAll.Range("F47").Value = 3727
